# This workbook is a price-tracking sheet: column A holds the product
# reference, columns B..BO hold one historical price snapshot per scrape
# (each header cell is the scrape timestamp), and the sheet ends with a
# "nom" (product name) column followed by a "url_produit" column.
#
# A new scrape happened, so a new snapshot column needs to be inserted
# right before "nom": it is stamped with the new timestamp in the header
# and, for every product row, simply repeats the most recent known price
# (i.e. a copy of the last existing price column, BO).
#
# Concretely this is an "insert column before BP" operation, followed by
# filling the freshly inserted BP column with the same values as BO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-30 18:25:56"

$lastPriceCol = 67   # column BO: the most recent existing price snapshot
$newCol       = 68   # column BP: where the new snapshot column goes
$lastRow      = 206  # last data row (dimension was A1:BQ206)

# Insert a new blank column at BP; this shifts "nom" (was BP) to BQ and
# "url_produit" (was BQ) to BR automatically, along with all their data.
$ws.Columns("BP").Insert()

# Stamp the header of the newly inserted column with the new scrape time.
$ws.Cells.Item(1, $newCol).Value2 = $newTimestamp

# For every data row, copy the latest known price (column BO) into the
# newly inserted column (BP) so the new snapshot reflects the last price.
for ($r = 2; $r -le $lastRow; $r++) {
    $lastPrice = $ws.Cells.Item($r, $lastPriceCol).Value2
    if ($lastPrice -ne "") {
        $ws.Cells.Item($r, $newCol).Value2 = $lastPrice
    }
}
